$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '26.269.27'
Set-TextValue 2 5 '  +0.78%  '
Set-TextValue 3 4 '1.677.70'
Set-TextValue 3 5 '  +0.56%  '
Set-TextValue 5 4 '217.43'
Set-TextValue 5 5 '  +0.23%  '
Set-TextValue 6 4 '0.5335'
Set-TextValue 6 5 '  +4.41%  '
Set-TextValue 7 4 '1.007'
Set-TextValue 7 5 '  +0.16%  '
Set-TextValue 8 5 '  +0.98%  '
Set-TextValue 9 4 '0.06466'
Set-TextValue 9 5 '  +0.82%  '
Set-TextValue 10 4 '21.89'
Set-TextValue 10 5 '  -0.10%  '
Set-TextValue 11 4 '0.07533'
Set-TextValue 11 5 '  +1.36%  '
Set-TextValue 12 4 '1.679.32'
Set-TextValue 12 5 '  +0.45%  '
Set-TextValue 13 5 '  +0.38%  '
Set-TextValue 14 4 '0.5765'
Set-TextValue 14 5 '  -1.70%  '
Set-TextValue 15 4 '0.000008454'
Set-TextValue 15 5 '  -1.32%  '
Set-TextValue 16 4 '64.67'
Set-TextValue 16 5 '  +0.53%  '
Set-TextValue 17 4 '26.293.48'
Set-TextValue 17 5 '  +0.81%  '
Set-TextValue 18 4 '4.899'
Set-TextValue 18 5 '  -0.87%  '
Set-TextValue 19 5 '  +0.24%  '
Set-TextValue 20 5 '  +0.89%  '
Set-TextValue 21 4 '189.91'
Set-TextValue 21 5 '  -0.39%  '
Set-TextValue 22 4 '6.199'
Set-TextValue 22 5 '  -0.47%  '
Set-TextValue 23 4 '1.007'
Set-TextValue 23 5 '  +0.11%  '
Set-TextValue 24 4 '145.59'
Set-TextValue 24 5 '  +0.24%  '
Set-TextValue 25 4 '7.820'
Set-TextValue 25 5 '  +2.75%  '
Set-TextValue 26 4 '0.1268'
Set-TextValue 26 5 '  +5.70%  '
Set-TextValue 27 4 '15.74'
Set-TextValue 27 5 '  +0.70%  '
Set-TextValue 28 4 '0.06469'
Set-TextValue 28 5 '  -3.12%  '
Set-TextValue 29 4 '1.386'
Set-TextValue 29 5 '  +5.34%  '
Set-TextValue 30 4 '1.318'
Set-TextValue 30 5 '  +0.23%  '
Set-TextValue 31 4 '3.581'
Set-TextValue 31 5 '  +1.25%  '
Set-TextValue 32 4 '3.589'
Set-TextValue 32 5 '  +2.05%  '
Set-TextValue 33 4 '1.658'
Set-TextValue 33 5 '  +0.57%  '
Set-TextValue 34 5 '  +1.30%  '
Set-TextValue 35 4 '0.6182'
Set-TextValue 35 5 '  +1.34%  '
Set-TextValue 36 4 '2.401'
Set-TextValue 36 5 '  +1.38%  '
Set-TextValue 37 4 '2.715'
Set-TextValue 37 5 '  +0.02%  '
Set-TextValue 38 4 '6.238'
Set-TextValue 38 5 '  +0.24%  '
Set-TextValue 39 4 '1.111.18'
Set-TextValue 39 5 '  +2.32%  '
Set-TextValue 40 4 '0.01622'
Set-TextValue 40 5 '  +1.24%  '
Set-TextValue 41 4 '0.8710'
Set-TextValue 41 5 '  +0.72%  '
Set-TextValue 42 4 '1.014'
Set-TextValue 42 5 '  +0.50%  '
Set-TextValue 43 5 '  -0.39%  '
Set-TextValue 44 4 '1.827.74'
Set-TextValue 44 5 '  +0.62%  '
Set-TextValue 45 4 '0.00000000108'
Set-TextValue 45 5 '  -6.26%  '
Set-TextValue 46 4 '57.07'
Set-TextValue 46 5 '  +1.41%  '
Set-TextValue 47 4 '8.152'
Set-TextValue 47 5 '  +1.00%  '
Set-TextValue 48 4 '1.003'
Set-TextValue 48 5 '  -0.53%  '
Set-TextValue 49 5 '  +0.43%  '
Set-TextValue 50 2 'Mantle'
Set-TextValue 50 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 50 4 '0.4289'
Set-TextValue 50 5 '  +0.04%  '
Set-TextValue 51 2 'Aptos'
Set-TextValue 51 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 51 4 '6.077'
Set-TextValue 51 5 '  +0.56%  '
